$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns retain their literal text representation
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "42.892.53"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3
$ws.Range("D3").Value = "2.548.45"
$ws.Range("E3").Value = "  -0.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "303.64"
$ws.Range("E5").Value = "  +1.43%  "

# Row 6
$ws.Range("D6").Value = "97.92"
$ws.Range("E6").Value = "  +5.85%  "

# Row 7
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  +0.60%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "0.544"
$ws.Range("E9").Value = "  -0.77%  "

# Row 10
$ws.Range("D10").Value = "37.01"
$ws.Range("E10").Value = "  +3.39%  "

# Row 11
$ws.Range("D11").Value = "0.0826"
$ws.Range("E11").Value = "  +2.59%  "

# Row 12
$ws.Range("E12").Value = "  +4.10%  "

# Row 13
$ws.Range("D13").Value = "7.76"
$ws.Range("E13").Value = "  +1.25%  "

# Row 14
$ws.Range("D14").Value = "2.941.25"
$ws.Range("E14").Value = "  +0.12%  "

# Row 15
$ws.Range("D15").Value = "2.562.50"
$ws.Range("E15").Value = "  +0.71%  "

# Row 16
$ws.Range("D16").Value = "15.14"
$ws.Range("E16").Value = "  +7.38%  "

# Row 17
$ws.Range("D17").Value = "0.876"
$ws.Range("E17").Value = "  +0.70%  "

# Row 18
$ws.Range("D18").Value = "42.910.38"
$ws.Range("E18").Value = "  -0.36%  "

# Row 19
$ws.Range("D19").Value = "13.85"
$ws.Range("E19").Value = "  +6.18%  "

# Row 20
$ws.Range("E20").Value = "  +1.23%  "

# Row 21
$ws.Range("D21").Value = "6.57"
$ws.Range("E21").Value = "  -0.44%  "

# Row 22
$ws.Range("D22").Value = "71.88"
$ws.Range("E22").Value = "  +0.25%  "

# Row 23
$ws.Range("D23").Value = "254.45"
$ws.Range("E23").Value = "  -0.77%  "

# Row 24
$ws.Range("D24").Value = "2.96"
$ws.Range("E24").Value = "  +1.86%  "

# Row 25
$ws.Range("E25").Value = "  -1.99%  "

# Row 26
$ws.Range("D26").Value = "28.01"
$ws.Range("E26").Value = "  -3.81%  "

# Row 27
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.06%  "

# Row 28
$ws.Range("D28").Value = "10.25"
$ws.Range("E28").Value = "  +2.19%  "

# Row 29
$ws.Range("D29").Value = "37.75"
$ws.Range("E29").Value = "  +1.31%  "

# Row 30
$ws.Range("E30").Value = "  -1.82%  "

# Row 31
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").Value = "  +3.56%  "

# Row 32
$ws.Range("D32").Value = "158.68"
$ws.Range("E32").Value = "  +3.71%  "

# Row 33
$ws.Range("D33").Value = "19.61"
$ws.Range("E33").Value = "  +15.12%  "

# Row 34
$ws.Range("E34").Value = "  -0.74%  "

# Row 35
$ws.Range("D35").Value = "0.0802"
$ws.Range("E35").Value = "  +0.61%  "

# Row 36
$ws.Range("D36").Value = "3.30"
$ws.Range("E36").Value = "  -2.22%  "

# Row 37
$ws.Range("E37").Value = "  -4.42%  "

# Row 38
$ws.Range("E38").Value = "  +1.57%  "

# Row 39
$ws.Range("D39").Value = "25.48"
$ws.Range("E39").Value = "  +9.57%  "

# Row 40
$ws.Range("E40").Value = "  -0.31%  "

# Row 41
$ws.Range("D41").Value = "2.11"
$ws.Range("E41").Value = "  +32.78%  "

# Row 42
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").Value = "3.42"
$ws.Range("E42").Value = "  -0.36%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "3.89"
$ws.Range("E43").Value = "  -0.26%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.094.61"
$ws.Range("E44").Value = "  +0.77%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0306"
$ws.Range("E45").Value = "  -1.48%  "

# Row 46
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.10%  "

# Row 47
$ws.Range("D47").Value = "86.50"
$ws.Range("E47").Value = "  +2.33%  "

# Row 48
$ws.Range("D48").Value = "8.96"
$ws.Range("E48").Value = "  +0.66%  "

# Row 49
$ws.Range("D49").Value = "75.31"
$ws.Range("E49").Value = "  +9.34%  "

# Row 50
$ws.Range("D50").Value = "2.798.33"
$ws.Range("E50").Value = "  +0.11%  "

# Row 51
$ws.Range("D51").Value = "103.25"
$ws.Range("E51").Value = "  -1.39%  "
